# Update "Horarios actualizados Línea 141 - 554"
# This script updates the schedule data on sheet "LP1912" (Sheet1),
# refreshes the "Última actualización" timestamp on all three sheets,
# updates the "Total filas" count on Sheet1, and removes the last
# data row (row 13) which no longer exists in the refreshed feed.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "03:49:28"
$newTimestamp = "04:01:12"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTimestamp"
$ws1.Range("A3").Value = "Total filas: 7"

# New data values for rows 6-12 (row 13 removed)
$data = @(
    @("04:02", "81_EL PELIGRO", 1),
    @("04:48", "81_EL PELIGRO", 47),
    @("04:53", "11_ETCHEVERRY", 52),
    @("05:17", "17_ROMERO", 76),
    @("05:22", "23_HERNANDEZ", 81),
    @("05:46", "14_ABASTO", 105),
    @("05:47", "17_ROMERO", 106)
)

$row = 6
foreach ($entry in $data) {
    $ws1.Cells.Item($row, 1).Value = $newTimestamp
    $ws1.Cells.Item($row, 2).Value = $entry[0]
    $ws1.Cells.Item($row, 3).Value = $entry[1]
    $ws1.Cells.Item($row, 4).Value = $entry[2]
    $ws1.Cells.Item($row, 5).Value = "LP1912"
    $row = $row + 1
}

# Remove old row 13 (data no longer present after refresh)
$ws1.Rows.Item(13).Delete()

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTimestamp"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTimestamp"
